# Add two new test case rows to the ACHData sheet of the BWP Bootstrap
# NormalizedSharedData workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ACHData")

# Row 6: ID 5 - Corporate Check Account Mismatch scenario
$ws.Range("A6").Value = "5"
$ws.Range("B6").Value = "95125489"
$ws.Range("D6").Value = "256072691"
$ws.Range("F6").Value = "999999999"
$ws.Range("G6").Value = "Corporate Check Acc Mismatch"
$ws.Range("C6").Value = "95125480"

# Row 7: ID 6 - Personal Checking Account Mismatch scenario
$ws.Range("A7").Value = "6"
$ws.Range("B7").Value = "95125489"
$ws.Range("D7").Value = "256072691"
$ws.Range("E7").Value = "1"
$ws.Range("G7").Value = "Personal Checking Acc Mismatch"
$ws.Range("C7").Value = "95125480"

# Match row height used in the target workbook for the new wrapped rows
$ws.Rows.Item(6).RowHeight = 28.8
$ws.Rows.Item(7).RowHeight = 28.8

# Update selection / active cell to reflect the last edited cell
$ws.Range("G7").Select()
